$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "309.47"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.12%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.03"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.44%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.122"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.32%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.57%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.402"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.10%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.305"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.62%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.868"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.93%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.983"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "5.50%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9254"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.13%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1106"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-8.13%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1876"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.38%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08821"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.66%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03295"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.28%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09574"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.59%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001393"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.98%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006208"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.392"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.04%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.09%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.374"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "20.38%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04346"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.79%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.75%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004274"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.33%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001403"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "7.95%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02138"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.14%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04973"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.64%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007575"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.36%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1353"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.23%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008506"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.71%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.21%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.007988"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-7.09%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006590"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.56%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.22%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "14.16%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "20.52%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.22%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.22%"
